$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the list of supermarkets / factors.
# Row 2 (Supermercado / 8) stays the same.
$ws.Range("A3").Value = "Hipermercado Metro Independencia"
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = "Plaza Vea Izaguirre"
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = "Tottus Mega Plaza"
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = "Plaza Vea Los Olivos"
$ws.Range("B6").Value = 1

$ws.Range("A7").Value = "Tottus Los Olivos"
$ws.Range("B7").Value = 1

$ws.Range("A8").Value = "Makro Plaza Lima Norte"
$ws.Range("B8").Value = 1

$ws.Range("A9").Value = "Makro Comas"
$ws.Range("B9").Value = 1
